# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
#
# Most cells are simple literal replacements. A handful of Price (column D)
# values are plain decimals ("7.00", "1.00", ...) whose trailing zero Excel
# would otherwise silently drop by auto-converting the cell to a Number on
# assignment. For those we use the standard leading-apostrophe quote-prefix
# to force text storage (matching the source data, which is text), then
# reset the cell's Style back to Normal so no stray "quote prefix" cell
# formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.071.77'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '3.771.69'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '''629.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '''165.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('D7').Value = '3.770.10'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('E10').Value = '  -2.43%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '''6.81'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.22%  '
$ws.Range('E13').Value = '  -5.14%  '
$ws.Range('D14').Value = '''34.84'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.82%  '
$ws.Range('D15').Value = '4.405.77'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').Value = '3.774.70'
$ws.Range('E16').Value = '  +5.00%  '
$ws.Range('D17').Value = '69.062.56'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '''17.66'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.76%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').Value = '''7.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.26%  '
$ws.Range('D21').Value = '''468.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').Value = '''9.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').Value = '''0.703'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('D24').Value = '''82.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.25%  '
$ws.Range('E25').Value = '  -8.26%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('D28').Value = '''10.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').Value = '3.920.09'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('D31').Value = '''2.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.38%  '
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('E33').Value = '  -3.18%  '
$ws.Range('E34').Value = '  +19.07%  '
$ws.Range('D35').Value = '''28.41'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = '3.724.07'
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('D38').Value = '''8.87'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.95%  '
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').Value = '''3.25'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('E41').Value = '  -2.70%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '''0.963'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.88%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '''1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E45').Value = '  +4.27%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = '''156.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = '''43.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '''46.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D49').Value = '''1.41'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.66%  '
$ws.Range('E50').Value = '  -2.43%  '
$ws.Range('D51').Value = '''8.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.54%  '
